# Update Leave Card 8/3/2023 4:40 PM
# Adds two new SICK LEAVE (SL(1-0-0)) entries for REOSA, CECILIA ANAY's leave
# card on Sheet1: one on the existing 7/1/2023 table row, and a brand-new
# table row (inserted immediately below it) for the second entry -
# both returned to work 1 day later on 7/11/2023 and 7/21/2023 respectively.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$lo = $ws.ListObjects.Item(1)

# --- Insert a brand-new table row at sheet row 112 -------------------------
# This pushes the existing row 112 (and everything below it, through the
# table's final helper row) down by one, growing the table from A8:K155 to
# A8:K156, exactly like using the "Insert Table Rows Above" command in Excel.
$ws.Rows.Item(112).Insert()
$lo.Resize($ws.Range("A8:K156"))

# The newly blank row 112 needs the same cell formatting/styling as the
# other plain data rows in the table (e.g. row 111 immediately above it).
$ws.Range("A111:K111").Copy()
$ws.Range("A112:K112").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# Restore the calculated "EARNED " column formula on the new row 112, and
# re-qualify it on the table's trailing helper row (156) which the resize
# above can leave pointing at an ambiguous unqualified column reference.
$ws.Range("G112").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'
$ws.Range("G156").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# --- Row 111 (PERIOD 7/1/2023, already present): log SL(1-0-0), 1 day ------
$ws.Range("B111").Value = "SL(1-0-0)"
$ws.Range("H111").Value = 1
$ws.Range("K110").Copy()
$ws.Range("K111").PasteSpecial(-4122)        # reuse the date-formatted style
$excel.CutCopyMode = 0
$ws.Range("K111").Value = 45118              # 7/11/2023 - date reported back

# --- Row 112 (new row, PERIOD left blank): log SL(1-0-0), 1 day ------------
$ws.Range("B112").Value = "SL(1-0-0)"
$ws.Range("H112").Value = 1
$ws.Range("K110").Copy()
$ws.Range("K112").PasteSpecial(-4122)        # reuse the date-formatted style
$excel.CutCopyMode = 0
$ws.Range("K112").Value = 45128              # 7/21/2023 - date reported back

# Reflect where the author's cursor ended up after making the edit.
$ws.Range("K112").Select()

$excel.Calculate()
